$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'item_link'
$ws.Range('B1').Value = 'img_path'
$ws.Range('C1').Value = 'best_review'
$ws.Range('D1').Value = 'title'
$ws.Range('A2').Value = 'https://link.coupang.com/a/cmgjus'
$ws.Range('B2').Value = 'http://thumbnail7.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/065f/32cfe83ec2d5f60dd258b5d5e109abc004e5ebaafe3e72064ff23e7ab855.jpg'
$ws.Range('C2').Value = '빠른 배송 감사합니다
인터넷 쇼핑으로 의류를 구매하고 실망한 경험이 있었는데 만족합니다
사진과 같은 품질입니다
원래 95를 입는데 100 사이즈가 저한테는 잘 맞습니다
가볍고 몸에 잘 맞고 겨드랑이부터 옆구리까지 매쉬 원단으로 통풍이 좋습니다
마감도 좋고 제품의 완성도도 좋아요
주머니도 깊어서 실용적이예요
방수기능은 아직 확인치 못했고 스판 기능은 없습니다
지퍼도 부드럽고 가벼운 운동이나 산책할때 입기에는 실망시키지안는 제품입니다
강추합니다
많이 파시고 대박나세요'
$ws.Range('D2').Value = '모렉스 남성 포르테 바람막이 모자 분리형 등산복 잠바 작업복 점퍼 자켓, 239_회색, 100'
$ws.Range('A3').Value = 'https://link.coupang.com/a/cmgjvF'
$ws.Range('B3').Value = 'http://thumbnail10.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/f946/97b84ec7d4b48c3b6fa58a9dbf36fb7c41c152cf18eada1c4f565c6248c1.jpg'
$ws.Range('C3').Value = '사용 후 솔직히 느낀대로 남겨봅니다!
구매 시 도움 되셨으면 좋겠습니다^^
** 나중에 재구매할 때 보기 위해서 최대한 꼼꼼히 적어둡니다. 다른분들에게 구매 시 도움 되셨으면 하는 바램 입니다 **
주문일 24.11.14
배송일 24.11.14
●구매당시금액 2 8 2 0 0
*구매사유
아직 겨울이 오기전이라 기모 들어간 옷만 너무 많은거 같아서 편하게 막 입을 옷을 찾다가 이곳저곳 활용하기 좋아보여 바로 구매했습니다.
*장점
피부에 직접 닿기 때문에 제일 중요한 소재,디자인,활용도 순을 중점으로 보는데 우선 디자인이 마음에 들어 구매했습니다.
비닐포장으로 왔고 전체적으로 마감,오염등 없이 잘 도착했습니다!
까슬거리거나 상표 등 거슬리는것은 없었습니다.
재질은 만지면 바스락 거리는 소리가 나고 , 냄새는 나지 않았습니다!
다만 실밥들이 좀 붙어 있었는데 툭툭 털면 떨어지는 정도 였습니다!
목부분의 지퍼는 부드럽게 잘 올라오고 내려가서 상황에 따라서 조절하면 될 듯 합니다!
앞에 주머니는 꽤 깊어서 너무 무거운 물건을 넣으면 배가 부각되어보일거 같아서 간단히 카드,이어폰 등 넣기 좋을 듯 합니다.
가운데 하단에는 조일 수 있게 되어 있는데 제 체형에는 조이지 않고 그냥 입는게 보기에도 좋고 편할꺼 같습니다!
하의는 현재 70키로 정도 나가고 , 하체에 살이 많은 타입인데
어제 검정 L-XL로 구매 했었는데 날씬해 보여 만족해서
한사이즈 업해서 베이지 색상으로 또! 구매 했습니다.
L 사이즈와 비교했을 때 우선 상의는 확실히 넉넉했고 , 손목에 시보리도 남아서 손을 걷어봤을 때 널럴해서 잘 여며서 넘기는게 아니라면 스르륵 풀려서 내려올 정도 였습니다!
하의는 앉았다 일어났을 때 전~혀 불편함 없이 움직일 수 있어서 좋았습니다. 자주 손이 갈꺼 같습니다^^!
아직 세탁,건조기 돌려보지는 않았는데 혹시 건조기까지 사용해서 줄어드는 상황이 오면 상품평에 추가해두겠습니다~!
두께가 적당히 얇지도 두껍지도 않아서 날씨가 지금보다는 조금 더 추워지면 패딩조끼와 같이 잘 활용해서 입을 수 있을 듯 합니다!
*단점
하루차이로 한사이즈 업해서 다른색상 구매했는데 금액이 올랐습니다 ㅠㅠ
내돈내산으로 제 상품평이 도움이 되셨다면 도움돼요 한번 부탁드립니다~^^'
$ws.Range('D3').Value = '잘빠진 여성용 어반 아노락 반집업 조거 상하의 세트, S~M, 카키'
$ws.Range('A4').Value = 'https://link.coupang.com/a/cmgjww'
$ws.Range('B4').Value = 'http://thumbnail10.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/9b89/ad8bbbf576af4ccf68b42d755e0223480f32ca1b1c5535c862c042243e9f.JPG'
$ws.Range('C4').Value = '깔끔하고 아주편하게 입기 좋아요.
컬러도 너무마음에들고 전반적으로 디자인이 아주깔끔해요,

대충 아무 바지에 걸쳐입기도좋아서 동네마트가거나 캠핑갈때, 산책에도좋아요.
포켓도 많고 가벼워서 사용하기도 최고!

여성용은 허리에 라인이 잡혀 부해보이지도 않는게 딱이었습니다.
휘뚤마뚤로 완전완전 추천해요!!!

심지어가격도 너무나 착해서 부담없더라구요.'
$ws.Range('D4').Value = '여성 바람막이 아웃도어 등산복 자켓, 핑크, 2XL'
$ws.Range('A5').Value = 'https://link.coupang.com/a/cmgjxy'
$ws.Range('B5').Value = 'http://thumbnail8.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/ba8c/78a259c7231dc46f31f8e905dbed70e232a7641595ee1de0dbde944dfac6.jpg'
$ws.Range('C5').Value = '민트스쿨 남성용 아웃도어 등산복 바람막이, 진그레이 5XL를 구매했어요. 남편이 체구가 있는 편이라 5XL 사이즈를 선택했는데, 정말 잘 맞더라고요. 색상은 호불호가 없는 진한 그레이라서 어떤 옷과도 잘 어울리고, 남편도 마음에 들어 했어요.
이 바람막이는 소매 주머니 지퍼 부분이 정말 견고하게 만들어져 있어서 보통 저렴한 제품에서 느껴지는 싼 티가 나지 않아요. 주머니 지퍼가 튼튼하게 잘 작동하고, 디자인도 깔끔해서 실용성이 높아요. 남편이 자주 외출할 때 필요한 아이템이라 더욱 신경 써서 고른 것인데, 이렇게 만족스러운 제품을 찾게 되어 기쁘네요.
하지만 한 가지 고민이 생겼어요. 세탁 후 소매와 등판 부분에 주름이 잘 없어지지 않더라고요. 처음에는 세탁 후 자연스럽게 펴질 거라고 생각했는데, 기대와는 달리 주름이 남아있어서 조금 아쉬웠어요. 그래서 세탁소에 가서 한 번 다려달라고 할까 고민 중이에요. 남편은 그냥 입겠다고 하긴 했지만, 주름이 약간 신경 쓰이기는
전반적인 품질은 정말 좋았어요. 옷 전체적으로 보풀도 없고, 세탁 후 변형도 전혀 없어서 사용하기에 문제없어요. 일반적으로 아웃도어 의류는 세탁 후 변형이 생기기 쉬운데, 이 제품은 그런 점에서 아주 만족스러웠어요. 남편이 자주 입고 다니는데, 활동할 때도 편안하고 자유롭게 움직일 수 있어서 좋다고 하더라고요.
바람막이의 방풍 기능도 뛰어나서 바람이 많이 부는 날에도 편안하게 입을 수 있었어요. 남편이 등산이나 야외 활동을 할 때 유용하게 활용하고 있어서, 이 점도 큰 장점이에요. 디자인도 심플하고 세련돼서 일상적인 외출 시에도 잘 어울려서 자주 입고 나가고 있어요.
이번에 구입한 아웃도어 바람막이는 색상과 디자인, 기능성 모두 만족스러운 제품이에요. 주름 문제는 조금 아쉽지만, 전체적인 만족도는 높아서 잘 입고 다니고 있어요.'
$ws.Range('D5').Value = '민트스쿨 남성용 아웃도어 등산복 바람막이, 블랙, 3XL'
$ws.Range('A6').Value = 'https://link.coupang.com/a/cmgjy5'
$ws.Range('B6').Value = 'http://thumbnail7.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/169b/897ec8656a9c0c399c29cd115b71b0726531590881750add0905a15891c7.jpg'
$ws.Range('C6').Value = '잘 받았습니다.
퀄리티가 너무 좋아서 특별히 착용 소감을 적어봤습니다.
더 많은 친구들이 이 겨울 자켓을 구매할 수 있기를 바랍니다.
이 겨울 재킷을 처음 입었을 때 사방에 감싸인 포근함이 순간적으로 온몸에 훈훈한 온기가 느껴졌습니다.
저는 아웃도어 활동과 자동차 여행을 즐기는 사람으로서 실용적이면서도 편안한 외투가 여행에서 얼마나 중요한지 잘 알고 있습니다.
이 재킷은 외부에 바람막이 원단을 사용함으로써 방풍,방수 효과가 뛰어나며 그로인해 방한과보온성 또한 뛰어납니다.
이 재킷은 마치 제 개인 보호막처럼 빗물이 때려도 안으로 전혀 스며들지 않아 항상 건조하고 편안함을 유지시켜 줍니다.
눈내리는날 우산없이 눈내리는 풍경을 마음껏 즐길 수 있을겁니다.
그리고 그 방풍기능은 더욱 절 감탄하게 합니다.
찬바람이 부는 겨울날에도 따뜻한 온기가 나를 감싸고 있는 것을 느낄 수 있습니다.
마치 따뜻한 품에 안겨 있는 것 같습니다.
창문을 열고 운전할 때도 이 재킷은 역시나 훌륭합니다.
재킷의 커팅 디자인이 매우 인간적이며, 큰 사이즈의 선택이 내 체형에 완벽하게 맞아떨어져서 운전석에서 아무런 구속 없이 자유롭게 움직일 수 있게 해줍니다.
또한 가벼운 소재와 스크레치 방지 기능으로 장시간 주행 후에도 의류가 손상되거나 주행 안전에 지장을 줄 염려 없이 산뜻하고 편안함을 유지시켜 줍니다.
더욱 놀라운 것은 이 재킷의 디자인이 매우 스타일리시하다는 것입니다.
단순히 실용적인 아웃도어 아이템이 아니라 전체적인 이미지를 높일 수 있는 패션 아이템입니다.
캐주얼 팬츠나 청바지에 코디하기에도 부담 없이 다양한 장소를 소화할 수 있어 아웃도어 생활과 일상생활에서도 자신감과 멋을 유지할 수 있습니다.
이제 저는 거의 매일 이 옷을 입고 이 추운 겨울을 함께 할 것입니다.'
$ws.Range('D6').Value = '남자패딩 등산 코트 빅사이즈 바람막이 아웃도어자켓 방한 아우터 방수 방풍 보드복 겨울 스키복, 다크그레이 dark..., S'
$ws.Range('A7').Value = 'https://link.coupang.com/a/cmgjBD'
$ws.Range('B7').Value = 'http://thumbnail7.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/2a27/f75c9f3bf68a4c3fec9b6fc75d14d6d7bd11d247085aead2de23ec3a7489.jpg'
$ws.Range('C7').Value = '✨maomaodian바람막이 등산복 점퍼 남성 방수 아웃도어 봄 가을 탈부착 J1

✔️세 가지 컬러
블랙 그레이 화이트

✔️아빠의 키는 170kg, 몸무게는 70kg 정도입니다.

아빠의 배가 조금 더 크신 분들도 계시기 때문에 추천 사이즈에 참고하셔서 평소 아빠가 사신 셔츠 사이즈는 105로 해서 2xl로 선택했습니다. 착용해보니 잘 어울리는 것 같으니 참고하세요.

✔️이 바랑막이 받고 정말 놀랐습니다! 소재가 정말 짱짱해요!

정전기 방지
땀과 보온을 방지하다.
가벼운 소재
스타일리시함
다목적
데일리 캐주얼, 등산, 캠핑, 달리기, 스키.

✔️바람막이 소재는 생각보다 만졌을 때 부드러운데 구김이 잘 가는 소재도 아니에요.예전에 만져본 적이 있는 바람막이는 부드럽고 구겨져 보이기 쉬운 바람막이인데 다림질하기도 쉽지 않아요!그런데 이 바람막이 소재는 안 되고 바람막이는 전체적으로 평평해 보이고 만져도 불편하지 않아요.

입으셔도 깔끔해 보여요!!

아빠가 까만색을 좋아해서 평소에도 어두운 계열로 코디하기 때문에 이 색의 아빠도 좋아해요. 게다가 디자인도 깔끔하고 소재도 편해서 아버지가 굉장히 만족해하셨어요.

✔️검정색도 오래 입을 수 있어요.

아주 특이한 목 부분의 단추 디자인인데 단추가 채워져 있는 단추보다는 이 디자인이 특별하고 예쁘다고 생각해요.

✔️바람막이 레터링도 깔끔합니다. 빨래를 몇 번 하면 떨어질 것 같은 인쇄된 문자는 아니에요.

✔️포켓에 지퍼 디테일이 있어 물건을 넣어도 잃어버릴 염려가 없답니다. 소매도 여밈이 있는 디자인으로 바람막이 효과가 더욱 좋습니다.

✔️모자에 달린 스트링도 신축성 있는 디자인으로 깔끔해서 쉽게 망가지지 않아요.

✔️모자는 분리가 가능하므로 모자를 원하지 않으셔도 쉽게 지퍼를 통해 바로 분리하실 수 있습니다.

전체적으로 정말 만족스러웠어요! 이 가격에 그렇게 좋은 바람막이를 살 수 있다니 정말 가성비가 너무 좋다! 몇 벌 더 사야겠다!

내돈내산의 찐 상품후기 입니다

⭕️소중한 내 재산으로 구입해서 사용해보고 쓴 생생 후기 입니다
⭕️재 후기가 도움이 되었다면‘도움이 돼요‘를 꾸욱 눌러주세요.'
$ws.Range('D7').Value = 'maomaodian바람막이 등산복 점퍼 남성 방수 아웃도어 봄 가을 탈부착 J1, 블랙, L(95)'
$ws.Range('A8').Value = 'https://link.coupang.com/a/cmgjC7'
$ws.Range('B8').Value = 'http://thumbnail8.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/2448/4ab510406faa34aa4307bac3697ff868cd0cd59d78f8a7369b3b98b89e20.jpg'
$ws.Range('C8').Value = '남편이 저녁때 조금 쌀쌀할 때 입을만한 옷 필요하다고 구매해달라고 해서 구매하게 되었답니당~^^♡

검정색 회색 카키색 있는데 회색이 젤 무난해보이더라구요.
검정색은 많아서~^^ 화면에 제 사진색보다는 밝습니다.


✅️남편이 많이 마른 편이예요.
살이 하나도 없어서... 95 입었는데도 넉넉합니다.
안에도 옷 두벌 정도 입었는데도 헐렁한?^^ 여유있는?^^ 폼입니다.
옷이 크게 나온건 아니고 정사이즈 같은데.. 95도 넉넉하네요. 살이 없어서 그런듯 합니다.
살 조금 있으면 95사이즈 사면(정사이즈) 작을듯 해요.
저는 살 많은 편인 여잔데도 작음. ㅎㅎ (가슴부분.배부분 ㅎㅎ)

옷은 이쁘고 편하다고는 해요

안에 핸드폰이나 지갑 넣는 곳이 있어요 .
깊이도 괜찮고 넓이도 좋네요~~♡♡^^
여자옷은 이런게 잘 없는데.. 남자들은 가방을 안가지고 다녀서 그런지 옷에 주머니가 여기저기 많은 듯 해요.
밖에도 가슴쪽에 자크 주머니가 또 있구요.
손목에 찍찍이로 줄일 수도 있어요^^ 그걸로 소매 조절하는거일텐데..
왜 있는거죠?^^ 애도 아니고 ㅎㅎ
손목 부분을 꽉 조일 이유가 있나요?^^
현장에서 일하거나 추울 땐 필요할것도 같긴해요?

안에 보면 허리쪽도 크기를 조절할 수 있는 끈도 있더라구요?^^
이것도 바지도 아니고~~ 왜 있는거죠?^^

신기한것들이 많네요. ㅎㅎ
딸내미만 키워봐서 이런 여러기능이 있는 옷을 처음 보내요.
다 쓸 곳들이 있겠지요?^^

쿠팡에서 후기작성 하면서 여러가지 배우네요~^^
새로운것들을 알게 해주셔서 감사드려요~^^

저녁에 추울 때 입으려고 샀는데..
날씨가 갑자기 더워져서 못입고 있어요. 좀 더 얇은 바람막이 정도 사야할 날씨가 되어버려서..
안그래도 더 얇은걸로 사라고 하네요~^^♡

농에 들어갈 듯 합니다.
올 가을에 입어야겟어요.. ㅎㅎ

저는 이 옷이 맘에 듭니다~^^♡♡ 색깔도 이쁘고 주머드들도 많고 자크도 부드럽게 잘 올라가고 자크 주머니라서 물건 분실율도 적을것 같고 실밥 터진곳도 하나도 없고..

진짜 이쁜옷 잘 산듯 해요.'
$ws.Range('D8').Value = '맨하트 춘추 아웃도어 자켓 등산복 일상복 작업복 빅사이즈 MHWF-59, 100, 블랙'
$ws.Range('A9').Value = 'https://link.coupang.com/a/cmgjD7'
$ws.Range('B9').Value = 'http://thumbnail8.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/a91c/9e3d33d7cdb1dce0bae29598b4c2acd0a7f7ca170087fc6e59da70ef2751.jpg'
$ws.Range('C9').Value = '작업복으로 입으려고 구매했습니다. 가랭이부분,엉덩이 부분,허리 부분, 지퍼 부분이 특히 중요한데요 박음질이 엉망입니다. 같은 사이즈인데 바지밑단의 넓이와 길이가 다릅니다. 일하다 가랭이나 엉덩이부분 터질까봐 걱정이네요.
재질,디자인은 마음에 듭니다. 박음질 견고하고 깔끔하게 좀 만들어주세요.'
$ws.Range('D9').Value = '워킹맨 아울렛 남성 1+1 봄 가을 시즌 다용도 멀티 스판 캠핑바지 낚시복 작업복 근무복 캠핑복 B128, 카키색+카키색, 34'
$ws.Range('A10').Value = 'https://link.coupang.com/a/cmggzt'
$ws.Range('B10').Value = 'http://thumbnail8.coupangcdn.com/thumbnails/remote/292x292ex/image/vendor_inventory/a91c/9e3d33d7cdb1dce0bae29598b4c2acd0a7f7ca170087fc6e59da70ef2751.jpg'
$ws.Range('C10').Value = '작업복으로 입으려고 구매했습니다. 가랭이부분,엉덩이 부분,허리 부분, 지퍼 부분이 특히 중요한데요 박음질이 엉망입니다. 같은 사이즈인데 바지밑단의 넓이와 길이가 다릅니다. 일하다 가랭이나 엉덩이부분 터질까봐 걱정이네요.
재질,디자인은 마음에 듭니다. 박음질 견고하고 깔끔하게 좀 만들어주세요.'
$ws.Range('D10').Value = '워킹맨 아울렛 남성 1+1 봄 가을 시즌 다용도 멀티 스판 캠핑바지 낚시복 작업복 근무복 캠핑복 B128, 카키색+카키색, 34'
